$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# New JSON body text for F2: same as existing, but "name" value replaced
# with a literal "Gopi Appapuram" instead of the $RandomFullName token.
$newBody = @'
{
  "createdAt": "$RandomPastDate",
  "name": "Gopi Appapuram",
  "email": "$RandomEmail",
  "given_name": "$RandomFirstName",
  "last_ip": "$RandomComputerIP",
  "updated_at": "$RandomPastDate",
  "last_login": "$RandomFutureDate",
  "email_verified": "$RandomBooleanValue"
}
'@

$ws.Range("F2").Value = $newBody

# Update the view: scroll back to top-left A1 and move the selection to G2.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G2").Select()
